$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 83 (region R2, tech WDS_ST_EXISTING) is removed entirely; all rows
# below it (84-140) shift up by one row, which is exactly what
# EntireRow.Delete() with an upward shift does.
$ws.Rows(83).Delete()
